$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = 0.05750029602107444
$ws.Range("G1").Value = 0.07020772298987588
$ws.Range("H1").Value = 0.625818588736696
$ws.Range("I1").Value = 0.05756203460963144
$ws.Range("J1").Value = 1.6247340862779955
$ws.Range("K1").Value = 0.17303792644239602
$ws.Range("L1").Value = 3.6881204119295345
$ws.Range("B2").Value = 6.0
$ws.Range("F2").Value = 0.06702953528633723
$ws.Range("G2").Value = 0.3510386149493794
$ws.Range("I2").Value = 0.06710150550538956
$ws.Range("J2").Value = 8.123670431389979
$ws.Range("L2").Value = 18.440602059647674
$ws.Range("F3").Value = 0.17109315953540039
$ws.Range("G3").Value = 0.01755193074746897
$ws.Range("H3").Value = 0.312909294368348
$ws.Range("I3").Value = 0.17127686381020277
$ws.Range("J3").Value = 0.4061835215694989
$ws.Range("K3").Value = 0.08651896322119801
$ws.Range("L3").Value = 0.9220301029823836
$ws.Range("F4").Value = 0.08771231596435083
$ws.Range("G4").Value = 0.01755193074746897
$ws.Range("H4").Value = 0.312909294368348
$ws.Range("I4").Value = 0.08780649347231914
$ws.Range("J4").Value = 0.4061835215694989
$ws.Range("K4").Value = 0.08651896322119801
$ws.Range("L4").Value = 0.9220301029823836
$ws.Range("B5").Value = 3.0
$ws.Range("F5").Value = 1.1326800308482836
$ws.Range("G5").Value = 0.561661783919007
$ws.Range("H5").Value = 1.251637177473392
$ws.Range("I5").Value = 1.1338961996548866
$ws.Range("J5").Value = 12.997872690223964
$ws.Range("K5").Value = 0.34607585288479203
$ws.Range("L5").Value = 29.504963295436276
$ws.Range("B6").Value = 3.0
$ws.Range("F6").Value = 1.8950191720693081
$ws.Range("G6").Value = 0.561661783919007
$ws.Range("H6").Value = 3.1290929436834802
$ws.Range("I6").Value = 1.897053871315537
$ws.Range("J6").Value = 12.997872690223964
$ws.Range("K6").Value = 0.8651896322119801
$ws.Range("L6").Value = 29.504963295436276
$ws.Range("B7").Value = 5.0
$ws.Range("F7").Value = 3.205289571042944
$ws.Range("G7").Value = 0.42124633793925526
$ws.Range("H7").Value = 2.816183649315132
$ws.Range("I7").Value = 3.20873111948228
$ws.Range("J7").Value = 9.748404517667973
$ws.Range("K7").Value = 0.7786706689907821
$ws.Range("L7").Value = 22.128722471577206
$ws.Range("B8").Value = 0.0
$ws.Range("F8").Value = 5.381854448676835
$ws.Range("G8").Value = 0.42124633793925526
$ws.Range("H8").Value = 1.5645464718417401
$ws.Range("I8").Value = 5.3876329945361245
$ws.Range("J8").Value = 9.748404517667973
$ws.Range("K8").Value = 0.43259481610599004
$ws.Range("L8").Value = 22.128722471577206
$ws.Range("F9").Value = 0.4548046012966339
$ws.Range("G9").Value = 0.42124633793925526
$ws.Range("H9").Value = 2.503274354946784
$ws.Range("I9").Value = 0.4552929291157289
$ws.Range("J9").Value = 9.748404517667973
$ws.Range("K9").Value = 0.6921517057695841
$ws.Range("L9").Value = 22.128722471577206
$ws.Range("F10").Value = 5.089480062128999
$ws.Range("G10").Value = 0.5265579224240691
$ws.Range("H10").Value = 1.251637177473392
$ws.Range("I10").Value = 5.094944682961728
$ws.Range("J10").Value = 12.185505647084966
$ws.Range("K10").Value = 0.34607585288479203
$ws.Range("L10").Value = 27.66090308947151
$ws.Range("F11").Value = 0.974581288492787
$ws.Range("I11").Value = 0.9756277052479905
$ws.Range("B12").Value = 3.0
$ws.Range("F12").Value = 1.6243021474879784
$ws.Range("G12").Value = 0.14041544597975175
$ws.Range("H12").Value = 2.503274354946784
$ws.Range("I12").Value = 1.6260461754133173
$ws.Range("J12").Value = 3.249468172555991
$ws.Range("K12").Value = 0.6921517057695841
$ws.Range("L12").Value = 7.376240823859069
$ws.Range("B13").Value = 0.0
$ws.Range("F13").Value = 2.165736196650638
$ws.Range("G13").Value = 0.14041544597975175
$ws.Range("H13").Value = 3.7549115324201763
$ws.Range("I13").Value = 2.1680615672177566
$ws.Range("J13").Value = 3.249468172555991
$ws.Range("K13").Value = 1.038227558654376
$ws.Range("L13").Value = 7.376240823859069
$ws.Range("B14").Value = 2.0
$ws.Range("F14").Value = 1.2994417179903828
$ws.Range("G14").Value = 0.21062316896962763
$ws.Range("H14").Value = 2.503274354946784
$ws.Range("I14").Value = 1.300836940330654
$ws.Range("J14").Value = 4.874202258833987
$ws.Range("K14").Value = 0.6921517057695841
$ws.Range("L14").Value = 11.064361235788603
$ws.Range("F15").Value = 1.2994417179903828
$ws.Range("G15").Value = 0.3510386149493794
$ws.Range("H15").Value = 3.7549115324201763
$ws.Range("I15").Value = 1.300836940330654
$ws.Range("J15").Value = 8.123670431389979
$ws.Range("K15").Value = 1.038227558654376
$ws.Range("L15").Value = 18.440602059647674
$ws.Range("F16").Value = 2.382309816315702
$ws.Range("G16").Value = 0.561661783919007
$ws.Range("H16").Value = 6.2581858873669605
$ws.Range("I16").Value = 2.3848677239395326
$ws.Range("J16").Value = 12.997872690223964
$ws.Range("K16").Value = 1.7303792644239602
$ws.Range("L16").Value = 29.504963295436276
$ws.Range("F17").Value = 8.229797547272424
$ws.Range("G17").Value = 0.42124633793925526
$ws.Range("H17").Value = 5.006548709893568
$ws.Range("I17").Value = 8.238633955427474
$ws.Range("J17").Value = 9.748404517667973
$ws.Range("K17").Value = 1.3843034115391681
$ws.Range("L17").Value = 22.128722471577206
